{"js": "const replacements = [\n  [\"22\u00d739=858\", \"79\u00d788=6952\"],\n  [\"28\u00d757=1596\", \"13\u00d790=1170\"],\n  [\"89\u00d786=7654\", \"66\u00d784=5544\"],\n  [\"81\u00d763=5103\", \"76\u00d770=5320\"],\n  [\"15\u00d789=1335\", \"95\u00d789=8455\"],\n  [\"34\u00d765=2210\", \"95\u00d715=1425\"],\n  [\"84\u00d755=4620\", \"56\u00d720=1120\"],\n  [\"51\u00d715=765\", \"67\u00d746=3082\"],\n  [\"21\u00d789=1869\", \"36\u00d780=2880\"],\n  [\"30\u00d744=1320\", \"89\u00d718=1602\"],\n  [\"60\u00d733=1980\", \"40\u00d741=1640\"],\n  [\"48\u00d754=2592\", \"53\u00d747=2491\"],\n  [\"31\u00d799=3069\", \"46\u00d788=4048\"],\n  [\"49\u00d761=2989\", \"91\u00d746=4186\"],\n  [\"28\u00d753=1484\", \"52\u00d717=884\"],\n  [\"85\u00d758=4930\", \"80\u00d785=6800\"],\n  [\"72\u00d794=6768\", \"26\u00d783=2158\"],\n  [\"76\u00d782=6232\", \"39\u00d786=3354\"],\n  [\"62\u00d743=2666\", \"19\u00d739=741\"],\n  [\"90\u00d791=8190\", \"81\u00d773=5913\"],\n  [\"97\u00d766=6402\", \"27\u00d719=513\"],\n  [\"60\u00d726=1560\", \"36\u00d738=1368\"],\n  [\"95\u00d748=4560\", \"82\u00d749=4018\"],\n  [\"84\u00d778=6552\", \"49\u00d722=1078\"],\n  [\"59\u00d774=4366\", \"55\u00d774=4070\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"22\u00d739=858\";  New = \"79\u00d788=6952\" },\n    @{ Old = \"28\u00d757=1596\"; New = \"13\u00d790=1170\" },\n    @{ Old = \"89\u00d786=7654\"; New = \"66\u00d784=5544\" },\n    @{ Old = \"81\u00d763=5103\"; New = \"76\u00d770=5320\" },\n    @{ Old = \"15\u00d789=1335\"; New = \"95\u00d789=8455\" },\n    @{ Old = \"34\u00d765=2210\"; New = \"95\u00d715=1425\" },\n    @{ Old = \"84\u00d755=4620\"; New = \"56\u00d720=1120\" },\n    @{ Old = \"51\u00d715=765\";  New = \"67\u00d746=3082\" },\n    @{ Old = \"21\u00d789=1869\"; New = \"36\u00d780=2880\" },\n    @{ Old = \"30\u00d744=1320\"; New = \"89\u00d718=1602\" },\n    @{ Old = \"60\u00d733=1980\"; New = \"40\u00d741=1640\" },\n    @{ Old = \"48\u00d754=2592\"; New = \"53\u00d747=2491\" },\n    @{ Old = \"31\u00d799=3069\"; New = \"46\u00d788=4048\" },\n    @{ Old = \"49\u00d761=2989\"; New = \"91\u00d746=4186\" },\n    @{ Old = \"28\u00d753=1484\"; New = \"52\u00d717=884\" },\n    @{ Old = \"85\u00d758=4930\"; New = \"80\u00d785=6800\" },\n    @{ Old = \"72\u00d794=6768\"; New = \"26\u00d783=2158\" },\n    @{ Old = \"76\u00d782=6232\"; New = \"39\u00d786=3354\" },\n    @{ Old = \"62\u00d743=2666\"; New = \"19\u00d739=741\" },\n    @{ Old = \"90\u00d791=8190\"; New = \"81\u00d773=5913\" },\n    @{ Old = \"97\u00d766=6402\"; New = \"27\u00d719=513\" },\n    @{ Old = \"60\u00d726=1560\"; New = \"36\u00d738=1368\" },\n    @{ Old = \"95\u00d748=4560\"; New = \"82\u00d749=4018\" },\n    @{ Old = \"84\u00d778=6552\"; New = \"49\u00d722=1078\" },\n    @{ Old = \"59\u00d774=4366\"; New = \"55\u00d774=4070\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Execute([ref]$r.Old, $false, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)\n}\n"}
